$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Change 1: split "Jul 2021, " into "Jul " + "2021, " runs, then
#     retext the first part to "May " (same rPr/formatting preserved) ---
$shp1 = $s.Shapes.Item(7)
$tr1 = $shp1.TextFrame.TextRange
$full1 = $tr1.Text
$idx1 = $full1.IndexOf("Jul 2021, ") + 1
$part1 = $tr1.Characters($idx1, 4)
# Forces the run to split in two (same formatting re-applied) so that
# "Jul " and "2021, " become distinct <a:r> runs.
$part1.Font.Size = $part1.Font.Size
$part1.Text = "May "

# --- Change 2: merge "-spark.memory.fraction" + "=" + "0.4" runs into
#     a single run with text "-spark.memory.fraction=0.4" ---
$shp2 = $s.Shapes.Item(16)
$tr2 = $shp2.TextFrame.TextRange
$full2 = $tr2.Text
$target2 = "-spark.memory.fraction=0.4"
$idx2 = $full2.IndexOf($target2) + 1
$sub2 = $tr2.Characters($idx2, $target2.Length)
$sub2.Text = $target2
